# Gasoline.xlsx -- "assignment 1 and class work"
#
# 1) On "Gasoline Sales Data": fill in the two moving-average forecast
#    columns (C = 2-period moving average, D = 3-period moving average),
#    leaving the cells that have no valid prior data as literal #N/A
#    errors (exactly what Excel shows if you drag the AVERAGE formula up
#    into rows that don't have enough preceding data).
# 2) Add a second worksheet ("Sheet1") with the same Week/Sales data plus
#    a simple exponential-smoothing forecast column.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Gasoline Sales Data: 2-period moving average (column C) ----------
$ws1.Range("C4").Style = "Normal"
$ws1.Range("C4").Value = "#N/A"

$ws1.Range("C5:C15").Style = "Normal"
$ws1.Range("C5:C15").Formula = "=AVERAGE(B2:B3)"

# --- Gasoline Sales Data: 3-period moving average (column D) ----------
$ws1.Range("D4:D5").Style = "Normal"
$ws1.Range("D4:D5").Value = "#N/A"

$ws1.Range("D6:D15").Style = "Normal"
$ws1.Range("D6:D15").Formula = "=AVERAGE(B2:B4)"

# --- view tweaks on the original sheet ---------------------------------
$ws1.Range("A1:B13").Select()
$excel.ActiveWindow.Zoom = 150

# --- add the second sheet, after "Gasoline Sales Data" -----------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet1"

# Bring over the Week / Sales table (values first, then formats so the
# new cells share the same style records as the source sheet).
$ws1.Range("A1:B13").Copy()
$ws2.Range("A1:B13").PasteSpecial(-4104)
$ws1.Range("A1:B13").Copy()
$ws2.Range("A1:B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Columns.Item(2).AutoFit()

# Simple exponential smoothing forecast in column C.
$ws2.Range("C3").Style = "Normal"
$ws2.Range("C3").Value = "#N/A"

$ws2.Range("C4").Style = "Normal"
$ws2.Range("C4").Formula = "=B2"

$ws2.Range("C5:C14").Style = "Normal"
$ws2.Range("C5:C14").Formula = "=0.8*B3+0.2*C4"

$ws2.Activate()
$ws2.Range("L21").Select()
